$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @("Neutrophils", "Slamf1", "Slamf1", "Neutrophils",
      2, 0.6666666666666666, 0.2022626666666667, 0.606788, 0.9444141808118884, 0.9444141808118884,
      2, 0.6666666666666666, 0.2022626666666667, 0.606788, 0.9444141808118884, 0.9444141808118884,
      0.04091018632711112, 0.368191676944, 0.8919181449185903, 0.8919181449185903),
    @("Neutrophils", "Slamf1", "Slamf1", "Resolving-Mac",
      2, 0.6666666666666666, 0.2022626666666667, 0.606788, 0.9444141808118884, 0.9444141808118884,
      1, 0.3333333333333333, 0.01190466666666667, 0.035714, 0.05558581918811148, 0.05558581918811148,
      0.002407869625777778, 0.021670826632, 0.05249603589329805, 0.05249603589329805),
    @("Resolving-Mac", "Slamf1", "Slamf1", "Neutrophils",
      1, 0.3333333333333333, 0.01190466666666667, 0.035714, 0.05558581918811148, 0.05558581918811148,
      2, 0.6666666666666666, 0.2022626666666667, 0.606788, 0.9444141808118884, 0.9444141808118884,
      0.002407869625777778, 0.021670826632, 0.05249603589329805, 0.05249603589329805),
    @("Resolving-Mac", "Slamf1", "Slamf1", "Resolving-Mac",
      1, 0.3333333333333333, 0.01190466666666667, 0.035714, 0.05558581918811148, 0.05558581918811148,
      1, 0.3333333333333333, 0.01190466666666667, 0.035714, 0.05558581918811148, 0.05558581918811148,
      0.0001417210884444445, 0.001275489796, 0.003089783294813422, 0.003089783294813422)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 2
    $row = $rows[$i]
    for ($c = 1; $c -le $row.Length; $c++) {
        $ws.Cells.Item($r, $c).Value = $row[$c - 1]
    }
}
